# "modify for tcc ipc"
# Adds a new "DEFINE_WAIT" code sample block to the wait_hc sheet, extends the
# bordered code-box styling down to cover the whole block, enables iterative
# calculation, and updates the active sheet / selections to match.

$wb = $excel.ActiveWorkbook

$wsWaitHc   = $wb.Worksheets.Item(1)   # wait_hc
$wsWaitSam  = $wb.Worksheets.Item(2)   # wait_sam
$wsComplHc  = $wb.Worksheets.Item(3)   # completion_hc
$wsComplSam = $wb.Worksheets.Item(4)   # completion_sam

# --- Enable iterative calculation (iterateDelta = 1E-4) -------------------
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# --- wait_hc: insert the new DEFINE_WAIT code block ------------------------
# 8 rows are inserted above the old row 19 ("#define DECLARE_WAITQUEUE..."):
# 7 new lines of code plus a trailing blank separator row.
$wsWaitHc.Rows("19:26").Insert() | Out-Null

$wsWaitHc.Range("C19").Value = "#define DEFINE_WAIT_FUNC(name, function)                    \"
$wsWaitHc.Range("C20").Value = "    wait_queue_t name = {                                   \"
$wsWaitHc.Range("C21").Value = "        .private      = current,                            \"
$wsWaitHc.Range("C22").Value = "        .func         = function,                           \"
$wsWaitHc.Range("C23").Value = "        .task_list    = LIST_HEAD_INIT((name).task_list),   \"
$wsWaitHc.Range("C24").Value = "    }"
$wsWaitHc.Range("C25").Value = "#define DEFINE_WAIT(name) DEFINE_WAIT_FUNC(name, autoremove_wake_function)"

# Extend the bordered "code box" look (as already used on completion_hc) over
# the whole C5:C44 block: middle style first, then fix up the top/bottom caps.
$wsComplHc.Range("C6").Copy() | Out-Null
$wsWaitHc.Range("C5:C44").PasteSpecial(-4122) | Out-Null

$wsComplHc.Range("C5").Copy() | Out-Null
$wsWaitHc.Range("C5").PasteSpecial(-4122) | Out-Null

$wsComplHc.Range("C25").Copy() | Out-Null
$wsWaitHc.Range("C44").PasteSpecial(-4122) | Out-Null

# --- Update selections on the other sheets ---------------------------------
$wsWaitSam.Range("M10").Select() | Out-Null
$wsComplSam.Range("S32").Select() | Out-Null

# Activate wait_hc last so it becomes the active/selected tab, with D25
# selected (matches the saved view state in the workbook).
$wsWaitHc.Activate() | Out-Null
$wsWaitHc.Range("D25").Select() | Out-Null
